# Added new script to RCC
# Adds a new test case row (RCC003) to the "Test Cases" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate the new row (row 3) with the new RCC003 test case data.
$ws.Range("A3").Value = "RCC003"
$ws.Range("B3").Value = "ABC"
$ws.Range("C3").Value = "Verifying the invitation information"
$ws.Range("D3").Value = "Y"

# Reflect the new selection left behind after entering the data.
$ws.Range("E3").Select()
